$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 213 (weekly update), pushing the existing
# rows 213 and 214 down to 215 and 216 respectively.
$ws.Rows("213:214").Insert()

# New row 213: updated/latest price observation for "Provincia de Curicó"
$ws.Cells.Item(213, 1).Value  = 6
$ws.Cells.Item(213, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(213, 3).Value  = "Metropolitana"
$ws.Cells.Item(213, 4).Value  = 44890
$ws.Cells.Item(213, 5).Value  = 13
$ws.Cells.Item(213, 6).Value  = "Fruta"
$ws.Cells.Item(213, 7).Value  = 100101
$ws.Cells.Item(213, 8).Value  = "Berries"
$ws.Cells.Item(213, 9).Value  = 100101004
$ws.Cells.Item(213, 10).Value = "Frambuesa"
$ws.Cells.Item(213, 11).Value = "Sin especificar"
$ws.Cells.Item(213, 12).Value = "Primera"
$ws.Cells.Item(213, 13).Value = 400
$ws.Cells.Item(213, 14).Value = 10000
$ws.Cells.Item(213, 15).Value = 10000
$ws.Cells.Item(213, 16).Value = 10000
$ws.Cells.Item(213, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(213, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(213, 19).Value = 5000
$ws.Cells.Item(213, 20).Value = 2

# New row 214: new price observation for "Región del Maule"
$ws.Cells.Item(214, 1).Value  = 6
$ws.Cells.Item(214, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(214, 3).Value  = "Metropolitana"
$ws.Cells.Item(214, 4).Value  = 44890
$ws.Cells.Item(214, 5).Value  = 13
$ws.Cells.Item(214, 6).Value  = "Fruta"
$ws.Cells.Item(214, 7).Value  = 100101
$ws.Cells.Item(214, 8).Value  = "Berries"
$ws.Cells.Item(214, 9).Value  = 100101004
$ws.Cells.Item(214, 10).Value = "Frambuesa"
$ws.Cells.Item(214, 11).Value = "Sin especificar"
$ws.Cells.Item(214, 12).Value = "Primera"
$ws.Cells.Item(214, 13).Value = 350
$ws.Cells.Item(214, 14).Value = 9000
$ws.Cells.Item(214, 15).Value = 9000
$ws.Cells.Item(214, 16).Value = 9000
$ws.Cells.Item(214, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(214, 18).Value = "Región del Maule"
$ws.Cells.Item(214, 19).Value = 4500
$ws.Cells.Item(214, 20).Value = 2
